# November-2021.xlsx edit: "10 years Finalization data"
#
# The author added a second worksheet ("Sheet1") that holds a finalized,
# self-contained copy of the daily data table (previously only living at
# A9:K39 on "Data Harian - Table"), re-numbered so the header starts at
# row 1 and the data runs through row 31. The new sheet becomes the
# active/visible tab; the original sheet's selection/scroll position also
# shifts down toward the data table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after the existing one.
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Copy the header + 30 days of data (A9:K39 on sheet1) to A1:K31 on the
# new sheet, values and formatting together, exactly like the original
# table.
$ws1.Range("A9:K39").Copy($ws2.Range("A1"))

# The new sheet has no explicit column widths, so with word-wrap on and
# Excel's narrower default column width the cells re-wrap to two lines;
# match the resulting row height for the data rows (header stays single
# line).
$ws2.Range("A2:K31").RowHeight = 28.8

# Restore/normalize sheet1's view (keep gridlines visible) and move the
# selection down onto the data table, like the saved file shows.
$ws1.Activate()
$excel.ActiveWindow.DisplayGridlines = $true
$ws1.Range("A9:K39").Select()

# Make the new sheet the active tab with the whole table selected.
$ws2.Activate()
$ws2.Range("A1:K31").Select()
